$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.993.66"
$ws.Range("E2").Value = "  -0.19%  "

$ws.Range("D3").Value = "1.677.32"
$ws.Range("E3").Value = "  +0.22%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.13"
$ws.Range("E5").Value = "  -0.48%  "

$ws.Range("E6").Value = "  +1.56%  "

$ws.Range("E7").Value = "  +0.10%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("E9").Value = "  +0.24%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.34"
$ws.Range("E10").Value = "  +1.03%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0888"
$ws.Range("E11").Value = "  -0.33%  "

$ws.Range("D12").Value = "1.913.30"
$ws.Range("E12").Value = "  +0.19%  "

$ws.Range("D13").Value = "1.661.86"
$ws.Range("E13").Value = "  -0.70%  "

$ws.Range("E14").Value = "  +0.13%  "

$ws.Range("E15").Value = "  +1.55%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.75"
$ws.Range("E16").Value = "  -0.11%  "

$ws.Range("D17").Value = "27.011.00"
$ws.Range("E17").Value = "  -0.16%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.18"
$ws.Range("E18").Value = "  +5.97%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "236.62"
$ws.Range("E19").Value = "  +0.73%  "

$ws.Range("E20").Value = "  -0.40%  "

$ws.Range("E21").Value = "  +0.07%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.45"
$ws.Range("E22").Value = "  -0.28%  "

$ws.Range("E23").Value = "  -0.67%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.18"
$ws.Range("E24").Value = "  -2.76%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.08"
$ws.Range("E25").Value = "  +0.44%  "

$ws.Range("E26").Value = "  +0.87%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.17"
$ws.Range("E27").Value = "  +1.66%  "

$ws.Range("E28").Value = "  -1.58%  "

$ws.Range("E29").Value = "  +0.07%  "

$ws.Range("E30").Value = "  -0.02%  "

$ws.Range("E31").Value = "  -0.47%  "

$ws.Range("E32").Value = "  -0.22%  "

$ws.Range("D33").Value = "1.483.84"
$ws.Range("E33").Value = "  +2.15%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.17"
$ws.Range("E34").Value = "  +1.06%  "

$ws.Range("E35").Value = "  +5.11%  "

$ws.Range("E36").Value = "  +0.24%  "

$ws.Range("E38").Value = "  +2.70%  "

$ws.Range("E39").Value = "  +1.31%  "

$ws.Range("E40").Value = "  -3.55%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.02"
$ws.Range("E41").Value = "  +1.01%  "

$ws.Range("E42").Value = "  +0.12%  "

$ws.Range("E43").Value = "  +1.96%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "67.48"
$ws.Range("E44").Value = "  +2.60%  "

$ws.Range("D45").Value = "1.821.79"
$ws.Range("E45").Value = "  +0.14%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.781"
$ws.Range("E46").Value = "  +0.02%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.38"
$ws.Range("E47").Value = "  +0.02%  "

$ws.Range("E48").Value = "  +0.44%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.52"
$ws.Range("E49").Value = "  -0.66%  "

$ws.Range("E50").Value = "  +1.35%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.75"
$ws.Range("E51").Value = "  +1.15%  "
